$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price (D) and Volume(1h) (E) columns store plain text values that
# often look numeric (e.g. "0.604", "1.00", "61.507.96"). Excel normally
# auto-converts such strings typed into a cell into real numbers, which
# would change the stored cell type. Force the whole data range to Text
# format first so our writes stay text, then restore the default style
# afterwards so untouched/changed cells end up without an explicit
# number-format style (matching the original plain-text cells).
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

# Apply the updated values cell by cell
$ws.Range("D2").Value = "61.507.96"
$ws.Range("E2").Value = "  -3.61%  "
$ws.Range("D3").Value = "2.483.25"
$ws.Range("E3").Value = "  -6.09%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "552.94"
$ws.Range("E5").Value = "  -4.75%  "
$ws.Range("D6").Value = "146.91"
$ws.Range("E6").Value = "  -5.62%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "0.604"
$ws.Range("E8").Value = "  -2.86%  "
$ws.Range("D9").Value = "2.481.97"
$ws.Range("E9").Value = "  -6.04%  "
$ws.Range("E10").Value = "  -8.59%  "
$ws.Range("E11").Value = "  -1.48%  "
$ws.Range("D12").Value = "5.43"
$ws.Range("E12").Value = "  -6.59%  "
$ws.Range("E13").Value = "  -6.27%  "
$ws.Range("D14").Value = "26.28"
$ws.Range("E14").Value = "  -7.79%  "
$ws.Range("D15").Value = "2.927.44"
$ws.Range("E15").Value = "  -6.08%  "
$ws.Range("D16").Value = "0.0000167"
$ws.Range("E16").Value = "  -8.83%  "
$ws.Range("D17").Value = "61.421.28"
$ws.Range("E17").Value = "  -3.73%  "
$ws.Range("D18").Value = "2.463.18"
$ws.Range("E18").Value = "  -6.61%  "
$ws.Range("D19").Value = "11.18"
$ws.Range("E19").Value = "  -8.08%  "
$ws.Range("D20").Value = "7.01"
$ws.Range("E20").Value = "  -8.43%  "
$ws.Range("E21").Value = "  -7.05%  "
$ws.Range("D22").Value = "322.63"
$ws.Range("E22").Value = "  -6.55%  "
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("D24").Value = "1.82"
$ws.Range("E24").Value = "  -4.58%  "
$ws.Range("D25").Value = "64.03"
$ws.Range("E25").Value = "  -5.83%  "
$ws.Range("E26").Value = "  -9.86%  "
$ws.Range("D27").Value = "2.607.85"
$ws.Range("E27").Value = "  -5.42%  "
$ws.Range("D28").Value = "1.52"
$ws.Range("E28").Value = "  -5.59%  "
$ws.Range("B29").Value = "Binance-PegBSC-USD"
$ws.Range("C29").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.25%  "
$ws.Range("B30").Value = "Bittensor"
$ws.Range("C30").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D30").Value = "538.11"
$ws.Range("E30").Value = "  -10.99%  "
$ws.Range("E31").Value = "  -9.80%  "
$ws.Range("D32").Value = "7.61"
$ws.Range("E32").Value = "  -5.96%  "
$ws.Range("E33").Value = "  -5.78%  "
$ws.Range("E34").Value = "  -7.68%  "
$ws.Range("E35").Value = "  -8.99%  "
$ws.Range("E36").Value = "  -10.62%  "
$ws.Range("D37").Value = "4.88"
$ws.Range("E37").Value = "  -10.27%  "
$ws.Range("E38").Value = "  -0.06%  "
$ws.Range("E39").Value = "  -5.36%  "
$ws.Range("D40").Value = "18.57"
$ws.Range("E40").Value = "  -6.21%  "
$ws.Range("D41").Value = "149.10"
$ws.Range("E41").Value = "  -1.62%  "
$ws.Range("E42").Value = "  -8.94%  "
$ws.Range("E43").Value = "  +0.04%  "
$ws.Range("D44").Value = "40.38"
$ws.Range("E44").Value = "  -3.59%  "
$ws.Range("D45").Value = "2.35"
$ws.Range("E45").Value = "  -8.50%  "
$ws.Range("D46").Value = "148.17"
$ws.Range("E46").Value = "  -7.91%  "
$ws.Range("D47").Value = "3.64"
$ws.Range("E47").Value = "  -7.18%  "
$ws.Range("D48").Value = "21.10"
$ws.Range("E48").Value = "  -14.30%  "
$ws.Range("D49").Value = "0.0539"
$ws.Range("E49").Value = "  -8.50%  "
$ws.Range("D50").Value = "0.597"
$ws.Range("E50").Value = "  -6.00%  "
$ws.Range("E51").Value = "  -4.95%  "

# Restore the default style on the whole range so no stray text-format
# style attribute is left behind on any cell.
$dataRange.Style = "Normal"
